$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the written-out grid so that row/column labels are 1-based instead
# of 0-based ("This is Row No 1 and col no 1" ... "This is Row No 10 and
# col no 10") while keeping the same A1:J10 physical layout.
for ($r = 1; $r -le 10; $r++) {
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value = "This is Row No $r and col no $c"
    }
}

# Reflect the updated selection state stored with the sheet (whole sheet
# selected, as captured in the saved view state).
[void]$ws.Cells.Select()
